# Refresh market-price derived figures (currentAveragePrice*, LevePrice*, LeveProfit*)
# for a batch of leve rows across multiple job sheets, per the scheduled market-data run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!40 - "Stuck in the Moment" (Horn Glue)
$ws.Cells.Item(40, 8).Value = 3848.25
$ws.Cells.Item(40, 9).Value = 3848.25
$ws.Cells.Item(40, 11).Value = 3848.25
$ws.Cells.Item(40, 13).Value = -3673.25

# ALC!116 - "Growing Up" (Growth Formula Kappa)
$ws.Cells.Item(116, 8).Value = 6874.75
$ws.Cells.Item(116, 10).Value = 6833.3335
$ws.Cells.Item(116, 12).Value = 6833.3335
$ws.Cells.Item(116, 14).Value = -13717.3335

# ALC!137 - "Cutting Edge of Culinary Quality" (Magnesia Whetstone)
$ws.Cells.Item(137, 8).Value = 2607.0588
$ws.Cells.Item(137, 9).Value = 2744.3572
$ws.Cells.Item(137, 10).Value = 1966.3334
$ws.Cells.Item(137, 11).Value = 8233.071599999999
$ws.Cells.Item(137, 12).Value = 5899.0002
$ws.Cells.Item(137, 13).Value = -5683.071599999999
$ws.Cells.Item(137, 14).Value = -10999.0002

$ws = $wb.Worksheets.Item("ARM")
# ARM!2 - "Ain't Got No Ingots" (Bronze Ingot)
$ws.Cells.Item(2, 8).Value = 1459.2424
$ws.Cells.Item(2, 9).Value = 691.8421
$ws.Cells.Item(2, 11).Value = 691.8421
$ws.Cells.Item(2, 13).Value = -578.8421

# ARM!74 - "As the Bolt Flies" (Titanium Nugget)
$ws.Cells.Item(74, 8).Value = 1452.6786
$ws.Cells.Item(74, 9).Value = 1286.3158
$ws.Cells.Item(74, 11).Value = 1286.3158
$ws.Cells.Item(74, 13).Value = -412.3158000000001

# ARM!77 - "Heavy Metal Banned (L)" (Titanium Nugget)
$ws.Cells.Item(77, 8).Value = 1452.6786
$ws.Cells.Item(77, 9).Value = 1286.3158
$ws.Cells.Item(77, 11).Value = 6431.579000000001
$ws.Cells.Item(77, 13).Value = -2063.579000000001

# ARM!88 - "The Mast Chance" (Adamantite Rivets)
$ws.Cells.Item(88, 8).Value = 1598.8334
$ws.Cells.Item(88, 9).Value = 1310.4445
$ws.Cells.Item(88, 10).Value = 1771.8667
$ws.Cells.Item(88, 11).Value = 1310.4445
$ws.Cells.Item(88, 12).Value = 1771.8667
$ws.Cells.Item(88, 13).Value = -904.4445000000001
$ws.Cells.Item(88, 14).Value = -2583.8667

# ARM!91 - "The Rose and the Riveter (L)" (Adamantite Rivets)
$ws.Cells.Item(91, 8).Value = 1598.8334
$ws.Cells.Item(91, 9).Value = 1310.4445
$ws.Cells.Item(91, 10).Value = 1771.8667
$ws.Cells.Item(91, 11).Value = 1310.4445
$ws.Cells.Item(91, 12).Value = 1771.8667
$ws.Cells.Item(91, 13).Value = 93.55549999999994
$ws.Cells.Item(91, 14).Value = -4579.8667

# ARM!102 - "Smells of Rich Tama-hagane" (Tama-hagane Ingot)
$ws.Cells.Item(102, 8).Value = 45650.266
$ws.Cells.Item(102, 9).Value = 16419.777
$ws.Cells.Item(102, 10).Value = 89496
$ws.Cells.Item(102, 11).Value = 16419.777
$ws.Cells.Item(102, 12).Value = 89496
$ws.Cells.Item(102, 13).Value = -14797.777
$ws.Cells.Item(102, 14).Value = -92740

# ARM!110 - "Scheduled Maintenance" (Deepgold Ingot)
$ws.Cells.Item(110, 8).Value = 4081.2307
$ws.Cells.Item(110, 9).Value = 4006.2222
$ws.Cells.Item(110, 10).Value = 4250
$ws.Cells.Item(110, 11).Value = 4006.2222
$ws.Cells.Item(110, 12).Value = 4250
$ws.Cells.Item(110, 13).Value = -1961.2222
$ws.Cells.Item(110, 14).Value = -8340

# ARM!116 - "No Scope" (Titanbronze Ingot)
$ws.Cells.Item(116, 8).Value = 1459.2424
$ws.Cells.Item(116, 9).Value = 691.8421
$ws.Cells.Item(116, 11).Value = 691.8421
$ws.Cells.Item(116, 13).Value = 1602.1579

$ws = $wb.Worksheets.Item("BSM")
# BSM!3 - "Hells Bells" (Bronze Ingot)
$ws.Cells.Item(3, 8).Value = 1459.2424
$ws.Cells.Item(3, 9).Value = 691.8421
$ws.Cells.Item(3, 11).Value = 691.8421
$ws.Cells.Item(3, 13).Value = -577.8421

$ws = $wb.Worksheets.Item("CRP")
# CRP!31 - "Wall Not Found" (Walnut Lumber)
$ws.Cells.Item(31, 8).Value = 1982.9524
$ws.Cells.Item(31, 9).Value = 2036.6
$ws.Cells.Item(31, 10).Value = 910
$ws.Cells.Item(31, 11).Value = 2036.6
$ws.Cells.Item(31, 12).Value = 910
$ws.Cells.Item(31, 13).Value = -1741.6
$ws.Cells.Item(31, 14).Value = -1500

# CRP!34 - "Armoires of the Rich and Famous" (Walnut Lumber)
$ws.Cells.Item(34, 8).Value = 1982.9524
$ws.Cells.Item(34, 9).Value = 2036.6
$ws.Cells.Item(34, 10).Value = 910
$ws.Cells.Item(34, 11).Value = 2036.6
$ws.Cells.Item(34, 12).Value = 910
$ws.Cells.Item(34, 13).Value = -1834.6
$ws.Cells.Item(34, 14).Value = -1314

# CRP!58 - "You Do the Heavy Lifting" (Mahogany Lumber)
$ws.Cells.Item(58, 8).Value = 4531.7036
$ws.Cells.Item(58, 9).Value = 2705.2
$ws.Cells.Item(58, 10).Value = 5606.1177
$ws.Cells.Item(58, 11).Value = 2705.2
$ws.Cells.Item(58, 12).Value = 5606.1177
$ws.Cells.Item(58, 13).Value = -2502.2
$ws.Cells.Item(58, 14).Value = -6012.1177

# CRP!62 - "Splinter in the Sewers" (Cedar Lumber)
$ws.Cells.Item(62, 8).Value = 5854.25
$ws.Cells.Item(62, 9).Value = 3806.5
$ws.Cells.Item(62, 10).Value = 11997.5
$ws.Cells.Item(62, 11).Value = 3806.5
$ws.Cells.Item(62, 12).Value = 11997.5
$ws.Cells.Item(62, 13).Value = -3182.5
$ws.Cells.Item(62, 14).Value = -13245.5

# CRP!65 - "The Lumber of Their Discontent (L)" (Cedar Lumber)
$ws.Cells.Item(65, 8).Value = 5854.25
$ws.Cells.Item(65, 9).Value = 3806.5
$ws.Cells.Item(65, 10).Value = 11997.5
$ws.Cells.Item(65, 11).Value = 19032.5
$ws.Cells.Item(65, 12).Value = 59987.5
$ws.Cells.Item(65, 13).Value = -15912.5
$ws.Cells.Item(65, 14).Value = -66227.5

# CRP!132 - "Hull Lotta Damage" (Ginseng Lumber)
$ws.Cells.Item(132, 8).Value = 2472.1177
$ws.Cells.Item(132, 9).Value = 2224.1333
$ws.Cells.Item(132, 11).Value = 6672.3999
$ws.Cells.Item(132, 13).Value = -4142.3999

# CRP!134 - "Wood You Be Quiet" (Ceiba Lumber)
$ws.Cells.Item(134, 8).Value = 2435.4443
$ws.Cells.Item(134, 9).Value = 2253.4
$ws.Cells.Item(134, 11).Value = 6760.200000000001
$ws.Cells.Item(134, 13).Value = -4225.200000000001

# CRP!136 - "Turali Quality" (Dark Mahogany Lumber)
$ws.Cells.Item(136, 8).Value = 4531.7036
$ws.Cells.Item(136, 9).Value = 2705.2
$ws.Cells.Item(136, 10).Value = 5606.1177
$ws.Cells.Item(136, 11).Value = 8115.599999999999
$ws.Cells.Item(136, 12).Value = 16818.3531
$ws.Cells.Item(136, 13).Value = -5565.599999999999
$ws.Cells.Item(136, 14).Value = -21918.3531

$ws = $wb.Worksheets.Item("GSM")
# GSM!70 - "Sky Is the Limit" (Mythrite Ingot)
$ws.Cells.Item(70, 8).Value = 21427.455
$ws.Cells.Item(70, 9).Value = 24023.535
$ws.Cells.Item(70, 11).Value = 24023.535
$ws.Cells.Item(70, 13).Value = -23753.535

# GSM!73 - "Hulls of Broken Dreams (L)" (Mythrite Ingot)
$ws.Cells.Item(73, 8).Value = 21427.455
$ws.Cells.Item(73, 9).Value = 24023.535
$ws.Cells.Item(73, 11).Value = 24023.535
$ws.Cells.Item(73, 13).Value = -23087.535

# GSM!107 - "Whetstones for the Workers" (Hard Mudstone Whetstone)
$ws.Cells.Item(107, 8).Value = 1662.3889
$ws.Cells.Item(107, 9).Value = 1899.3077
$ws.Cells.Item(107, 10).Value = 1046.4
$ws.Cells.Item(107, 11).Value = 1899.3077
$ws.Cells.Item(107, 12).Value = 1046.4
$ws.Cells.Item(107, 13).Value = 20.69229999999993
$ws.Cells.Item(107, 14).Value = -4886.4

$ws = $wb.Worksheets.Item("LTW")
# LTW!40 - "Best Served Toad" (Toad Leather)
$ws.Cells.Item(40, 8).Value = 3864.8333
$ws.Cells.Item(40, 9).Value = 4035.4546
$ws.Cells.Item(40, 10).Value = 1988
$ws.Cells.Item(40, 11).Value = 4035.4546
$ws.Cells.Item(40, 12).Value = 1988
$ws.Cells.Item(40, 13).Value = -3899.4546
$ws.Cells.Item(40, 14).Value = -2260

# LTW!68 - "You Could Say It's a Moving Target" (Wyvern Leather)
$ws.Cells.Item(68, 8).Value = 2407
$ws.Cells.Item(68, 10).Value = 3090
$ws.Cells.Item(68, 12).Value = 3090
$ws.Cells.Item(68, 14).Value = -4588

# LTW!71 - "They Call It Bloody Mary (L)" (Wyvern Leather)
$ws.Cells.Item(71, 8).Value = 2407
$ws.Cells.Item(71, 10).Value = 3090
$ws.Cells.Item(71, 12).Value = 15450
$ws.Cells.Item(71, 14).Value = -22938

# LTW!82 - "Trainin' the Neck" (Dragon Leather)
$ws.Cells.Item(82, 8).Value = 1638
$ws.Cells.Item(82, 9).Value = 1357.6
$ws.Cells.Item(82, 10).Value = 1988.5
$ws.Cells.Item(82, 11).Value = 1357.6
$ws.Cells.Item(82, 12).Value = 1988.5
$ws.Cells.Item(82, 13).Value = -996.5999999999999
$ws.Cells.Item(82, 14).Value = -2710.5

# LTW!85 - "Training Is Only Skintight (L)" (Dragon Leather)
$ws.Cells.Item(85, 8).Value = 1638
$ws.Cells.Item(85, 9).Value = 1357.6
$ws.Cells.Item(85, 10).Value = 1988.5
$ws.Cells.Item(85, 11).Value = 1357.6
$ws.Cells.Item(85, 12).Value = 1988.5
$ws.Cells.Item(85, 13).Value = -109.5999999999999
$ws.Cells.Item(85, 14).Value = -4484.5

$ws = $wb.Worksheets.Item("WVR")
# WVR!80 - "Healing with Flair" (Hallowed Ramie Gaskins of Healing)
$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 14).ClearContents()

# WVR!83 - "Pants Fit for Battle (L)" (Hallowed Ramie Gaskins of Healing)
$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 14).ClearContents()

# WVR!107 - "Flax Wax" (Bright Linen Yarn)
$ws.Cells.Item(107, 8).Value = 906.2222
$ws.Cells.Item(107, 9).Value = 857.6667
$ws.Cells.Item(107, 11).Value = 2573.0001
$ws.Cells.Item(107, 13).Value = -653.0001000000002

# WVR!132 - "Comfy Cabins" (Snow Cotton Cloth)
$ws.Cells.Item(132, 8).Value = 1143.72
$ws.Cells.Item(132, 9).Value = 1075.4117
$ws.Cells.Item(132, 10).Value = 1288.875
$ws.Cells.Item(132, 11).Value = 3226.2351
$ws.Cells.Item(132, 12).Value = 3866.625
$ws.Cells.Item(132, 13).Value = -696.2351000000003
$ws.Cells.Item(132, 14).Value = -8926.625

# WVR!136 - "Weaving the Envelope" (Sarcenet Cloth)
$ws.Cells.Item(136, 8).Value = 5492.25
$ws.Cells.Item(136, 9).Value = 4991.143
$ws.Cells.Item(136, 11).Value = 14973.429
$ws.Cells.Item(136, 13).Value = -12423.429

Write-Host "Done updating leve profit figures."